# Duplicate the "4.3.2 (Jan 22nd)" worksheet to create a new
# "4.3.2 (Mar 2nd)" worksheet, placed right after the original and before
# the "4.3.2 (EMPTY)" sheet, then make it the active/selected sheet.

$wb = $excel.ActiveWorkbook

$sourceSheet = $wb.Worksheets.Item("4.3.2 (Jan 22nd)")

# Copy the source sheet to immediately after itself - this inserts the new
# sheet between "4.3.2 (Jan 22nd)" and "4.3.2 (EMPTY)".
$sourceSheet.Copy($null, $sourceSheet)

$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "4.3.2 (Mar 2nd)"

# Make the newly created sheet the active tab / selected cell, matching the
# author's last-saved UI state.
$newSheet.Activate()
$newSheet.Range("E28").Select()
